# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Balmung
# profit sheets. Only numeric value cells are touched; no formulas or
# formatting are involved since these columns are plain static snapshots.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2679.818
$ws.Range("J32").Value = 2919.889
$ws.Range("L32").Value = 2919.889
$ws.Range("N32").Value = -3571.889

$ws.Range("H80").Value = 1044.4546
$ws.Range("I80").Value = 874.75
$ws.Range("J80").Value = 1141.4286
$ws.Range("K80").Value = 2624.25
$ws.Range("L80").Value = 3424.2858
$ws.Range("M80").Value = -1626.25
$ws.Range("N80").Value = -5420.2858

$ws.Range("H83").Value = 1044.4546
$ws.Range("I83").Value = 874.75
$ws.Range("J83").Value = 1141.4286
$ws.Range("K83").Value = 7872.75
$ws.Range("L83").Value = 10272.8574
$ws.Range("M83").Value = -2880.75
$ws.Range("N83").Value = -20256.8574

$ws.Range("H131").Value = 7008.5
$ws.Range("I131").Value = 6264.1665
$ws.Range("K131").Value = 18792.4995
$ws.Range("M131").Value = -13752.4995

$ws.Range("H138").Value = 4346.528
$ws.Range("J138").Value = 2982.1455
$ws.Range("L138").Value = 8946.4365
$ws.Range("N138").Value = -19226.4365

$ws.Range("H141").Value = 1809.8889
$ws.Range("I141").Value = 1411.125
$ws.Range("K141").Value = 4233.375
$ws.Range("M141").Value = 946.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 185814.86
$ws.Range("I32").Value = 218489.27
$ws.Range("K32").Value = 218489.27
$ws.Range("M32").Value = -218202.27

$ws.Range("H97").Value = 8476.933999999999
$ws.Range("I97").Value = 8861.083000000001
$ws.Range("K97").Value = 8861.083000000001
$ws.Range("M97").Value = -8365.083000000001

$ws.Range("H122").Value = 1998
$ws.Range("I122").Value = 1998
$ws.Range("K122").Value = 5994
$ws.Range("M122").Value = -3544

$ws.Range("H132").Value = 2459.975
$ws.Range("I132").Value = 2131.625
$ws.Range("J132").Value = 3773.375
$ws.Range("K132").Value = 6394.875
$ws.Range("L132").Value = 11320.125
$ws.Range("M132").Value = -3864.875
$ws.Range("N132").Value = -16380.125

$ws.Range("H135").Value = 87486
$ws.Range("J135").Value = 87486
$ws.Range("L135").Value = 87486
$ws.Range("N135").Value = -97626

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2596.0667
$ws.Range("I86").Value = 1446.1111
$ws.Range("K86").Value = 1446.1111
$ws.Range("M86").Value = -323.1111000000001

$ws.Range("H89").Value = 2596.0667
$ws.Range("I89").Value = 1446.1111
$ws.Range("K89").Value = 7230.5555
$ws.Range("M89").Value = -1614.5555

$ws.Range("H105").Value = 18149.572
$ws.Range("I105").Value = 21109.6
$ws.Range("K105").Value = 21109.6
$ws.Range("M105").Value = -19362.6

$ws.Range("H134").Value = 3263.0952
$ws.Range("I134").Value = 3302.818
$ws.Range("J134").Value = 3219.4
$ws.Range("K134").Value = 9908.454000000002
$ws.Range("L134").Value = 9658.200000000001
$ws.Range("M134").Value = -7373.454000000002
$ws.Range("N134").Value = -14728.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2609.4082
$ws.Range("I31").Value = 3554.9167
$ws.Range("J31").Value = 2302.7568
$ws.Range("K31").Value = 3554.9167
$ws.Range("L31").Value = 2302.7568
$ws.Range("M31").Value = -3259.9167
$ws.Range("N31").Value = -2892.7568

$ws.Range("H34").Value = 2609.4082
$ws.Range("I34").Value = 3554.9167
$ws.Range("J34").Value = 2302.7568
$ws.Range("K34").Value = 3554.9167
$ws.Range("L34").Value = 2302.7568
$ws.Range("M34").Value = -3352.9167
$ws.Range("N34").Value = -2706.7568

$ws.Range("H122").Value = 1988.1621
$ws.Range("I122").Value = 1891.6
$ws.Range("K122").Value = 5674.799999999999
$ws.Range("M122").Value = -3224.799999999999

$ws.Range("H132").Value = 42333.88
$ws.Range("I132").Value = 54660.105
$ws.Range("K132").Value = 163980.315
$ws.Range("M132").Value = -161450.315

$ws.Range("H134").Value = 2231.375
$ws.Range("I134").Value = 1935.6364
$ws.Range("K134").Value = 5806.9092
$ws.Range("M134").Value = -3271.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3798.2222
$ws.Range("I34").Value = 294
$ws.Range("J34").Value = 4799.4287
$ws.Range("K34").Value = 882
$ws.Range("L34").Value = 14398.2861
$ws.Range("M34").Value = -798
$ws.Range("N34").Value = -14566.2861

$ws.Range("H86").Value = 1498.6
$ws.Range("I86").Value = 1499.3334
$ws.Range("J86").Value = 1497.5
$ws.Range("K86").Value = 4498.0002
$ws.Range("L86").Value = 4492.5
$ws.Range("M86").Value = -3312.0002
$ws.Range("N86").Value = -6864.5

$ws.Range("H89").Value = 1498.6
$ws.Range("I89").Value = 1499.3334
$ws.Range("J89").Value = 1497.5
$ws.Range("K89").Value = 13494.0006
$ws.Range("L89").Value = 13477.5
$ws.Range("M89").Value = -7566.000599999999
$ws.Range("N89").Value = -25333.5

$ws.Range("H113").Value = 247.16129
$ws.Range("I113").Value = 231.71428
$ws.Range("K113").Value = 695.14284
$ws.Range("M113").Value = 1474.85716

$ws.Range("H140").Value = 2159.0908
$ws.Range("I140").Value = 2075
$ws.Range("K140").Value = 6225
$ws.Range("M140").Value = -1045

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1974.5
$ws.Range("I113").Value = 1749.5
$ws.Range("K113").Value = 1749.5
$ws.Range("M113").Value = 420.5

$ws.Range("H132").Value = 1715.6364
$ws.Range("I132").Value = 1932.0769
$ws.Range("K132").Value = 5796.2307
$ws.Range("M132").Value = -3266.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1866.6666
$ws.Range("I68").Value = 1599
$ws.Range("J68").Value = 2000.5
$ws.Range("K68").Value = 1599
$ws.Range("L68").Value = 2000.5
$ws.Range("M68").Value = -850
$ws.Range("N68").Value = -3498.5

$ws.Range("H71").Value = 1866.6666
$ws.Range("I71").Value = 1599
$ws.Range("J71").Value = 2000.5
$ws.Range("K71").Value = 7995
$ws.Range("L71").Value = 10002.5
$ws.Range("M71").Value = -4251
$ws.Range("N71").Value = -17490.5

$ws.Range("H100").Value = 5309
$ws.Range("I100").Value = 4338.5
$ws.Range("K100").Value = 4338.5
$ws.Range("M100").Value = -3797.5

$ws.Range("H122").Value = 3171.375
$ws.Range("J122").Value = 3466.3333
$ws.Range("L122").Value = 10398.9999
$ws.Range("N122").Value = -15298.9999

$ws.Range("H132").Value = 2485.25
$ws.Range("J132").Value = 4664.6665
$ws.Range("L132").Value = 13993.9995
$ws.Range("N132").Value = -19053.9995

$ws.Range("H136").Value = 12734.689
$ws.Range("I136").Value = 40216.31
$ws.Range("K136").Value = 120648.93
$ws.Range("M136").Value = -118098.93

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 1500
$ws.Range("M51").Value = -990

$ws.Range("H113").Value = 402.95
$ws.Range("I113").Value = 398.29413
$ws.Range("K113").Value = 1194.88239
$ws.Range("M113").Value = 975.11761

$ws.Range("H132").Value = 1724.1794
$ws.Range("I132").Value = 1189.8572
$ws.Range("K132").Value = 3569.5716
$ws.Range("M132").Value = -1039.5716
